# Applies the "Updated cryptos list" data refresh to the crypto price table.
# Columns: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "54.197.10"
$ws.Range("E2").Value2 = "  -0.04%  "
$ws.Range("D3").Value2 = "2.262.85"
$ws.Range("E3").Value2 = "  -1.18%  "
$ws.Range("E4").Value2 = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "495.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +0.55%  "
$ws.Range("E7").Value2 = "  -0.20%  "
$ws.Range("E8").Value2 = "  -0.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0952"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  +0.26%  "
$ws.Range("E10").Value2 = "  +1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +3.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  +5.07%  "
$ws.Range("D14").Value2 = "2.663.25"
$ws.Range("E14").Value2 = "  -1.23%  "
$ws.Range("D15").Value2 = "54.177.74"
$ws.Range("E15").Value2 = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000129"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +0.13%  "
$ws.Range("D17").Value2 = "2.274.60"
$ws.Range("E17").Value2 = "  -1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "302.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +3.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  +1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +1.31%  "
$ws.Range("D30").Value2 = "0.0₃0687"
$ws.Range("E30").Value2 = "  -0.40%  "
$ws.Range("E31").Value2 = "  +0.60%  "
$ws.Range("E32").Value2 = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +0.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -0.10%  "
$ws.Range("E35").Value2 = "  +3.28%  "
$ws.Range("E36").Value2 = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.374"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +0.13%  "
$ws.Range("E39").Value2 = "  -0.41%  "
$ws.Range("B40").Value2 = "Filecoin"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -0.08%  "
$ws.Range("B41").Value2 = "RenderToken"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "124.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -2.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0492"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +1.39%  "
$ws.Range("E44").Value2 = "  +0.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.544"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "241.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +1.00%  "
$ws.Range("E47").Value2 = "  +0.19%  "
$ws.Range("E48").Value2 = "  +0.47%  "
$ws.Range("E49").Value2 = "  +0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -1.42%  "
$ws.Range("E51").Value2 = "  -0.41%  "
